$wb = $excel.ActiveWorkbook

# 1. Rename the "Include" sheets.
$wb.Worksheets.Item("Include from oBDS Operation K").Name = "Include #0"
$wb.Worksheets.Item("Include from ICD10GM").Name = "Include #1"

# 2. Update the Metadata sheet.
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" value (row 8, column B).
$ws.Cells.Item(8, 2).Value() = "2024-09-17T19:55:11+00:00"

# Shift rows 11-14 down to 12-15 to make room for the new "Jurisdiction" row
# (Description/Purpose/Copyright/Immutable all move down by one row).
for ($r = 14; $r -ge 11; $r--) {
    $ws.Cells.Item($r + 1, 1).Value() = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r + 1, 2).Value() = $ws.Cells.Item($r, 2).Value()
}

# Restore formatting on the newly extended row 15 (copy format from row 14,
# since the values-only copy above does not carry over cell formatting).
$fmtSrc = $ws.Range("A14:B14")
$fmtDst = $ws.Range("A15:B15")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)

# Write the new "Jurisdiction" row (row 11) with an empty value.
$ws.Cells.Item(11, 1).Value() = "Jurisdiction"
$ws.Cells.Item(11, 2).Value() = "'"

# Restore plain formatting on B11 (the leading apostrophe used above to force
# an explicit empty text value also marks the cell as quote-prefixed).
$fmtSrc2 = $ws.Range("B12")
$fmtSrc2.Copy()
$ws.Range("B11").PasteSpecial(-4122)
